# Auto-generated cell updates applying the commit diff to the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.053.16'
$ws.Range('E2').Value = '  +2.23%  '
$ws.Range('D3').Value = '3.809.74'
$ws.Range('E3').Value = '  +0.82%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Value = '''629.50'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.29%  '
$ws.Range('D6').Value = '''164.99'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.42%  '
$ws.Range('D7').Value = '3.806.84'
$ws.Range('E7').Value = '  +0.78%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +1.01%  '
$ws.Range('E10').Value = '  +3.01%  '
$ws.Range('E11').Value = '  +0.97%  '
$ws.Range('D12').Value = '''6.60'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.23%  '
$ws.Range('D13').Value = '''0.0000250'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('D14').Value = '''35.92'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.37%  '
$ws.Range('D15').Value = '4.447.93'
$ws.Range('E15').Value = '  +0.76%  '
$ws.Range('D16').Value = '3.931.05'
$ws.Range('E16').Value = '  +4.32%  '
$ws.Range('D17').Value = '69.000.70'
$ws.Range('E17').Value = '  +2.11%  '
$ws.Range('D18').Value = '''17.92'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.07%  '
$ws.Range('E19').Value = '  +1.47%  '
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('D21').Value = '''466.25'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.54%  '
$ws.Range('D22').Value = '''9.70'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').Value = '''0.707'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.12%  '
$ws.Range('E24').Value = '  +4.66%  '
$ws.Range('D25').Value = '''83.64'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('B26').Value = 'Fetch.AI'
$ws.Range('C26').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D26').Value = '''2.16'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.44%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').Value = '''11.97'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('D28').Value = '''10.08'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.77%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').Value = '3.958.77'
$ws.Range('E30').Value = '  +0.80%  '
$ws.Range('E31').Value = '  +4.17%  '
$ws.Range('D32').Value = '''2.22'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.53%  '
$ws.Range('D33').Value = '''7.27'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.02%  '
$ws.Range('D34').Value = '''29.16'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.51%  '
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('D36').Value = '''9.06'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.34%  '
$ws.Range('E37').Value = '  +3.73%  '
$ws.Range('E38').Value = '  +8.26%  '
$ws.Range('D39').Value = '''3.46'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.95%  '
$ws.Range('D40').Value = '''5.90'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.76%  '
$ws.Range('D41').Value = '''0.974'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.15%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').Value = '''157.46'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.64%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').Value = '''0.300'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.63%  '
$ws.Range('B46').Value = 'ONDO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D46').Value = '''1.42'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.09%  '
$ws.Range('B47').Value = 'Arweave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D47').Value = '''43.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.23%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').Value = '''46.86'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.10%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '''1.90'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.45%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '''8.43'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.47%  '
$ws.Range('E51').Value = '  +15.08%  '
